# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Updates columns I (DAMSLTag) and J (DialogAct) for the rows whose
# annotations changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{ Row = 7;   I = "b";  J = "Acknowledge (Backchannel)" },
    @{ Row = 14;  I = "ba"; J = "Appreciation" },
    @{ Row = 25;  I = "ba"; J = "Appreciation" },
    @{ Row = 26;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 35;  I = "ba"; J = "Appreciation" },
    @{ Row = 37;  I = "aa"; J = "Agree/Accept" },
    @{ Row = 38;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 41;  I = "%";  J = "Uninterpretable" },
    @{ Row = 44;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 55;  I = "aa"; J = "Agree/Accept" },
    @{ Row = 64;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 71;  I = "b";  J = "Acknowledge (Backchannel)" },
    @{ Row = 74;  I = "ba"; J = "Appreciation" },
    @{ Row = 87;  I = "ba"; J = "Appreciation" },
    @{ Row = 94;  I = "aa"; J = "Agree/Accept" },
    @{ Row = 95;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 107; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 108; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 111; I = "aa"; J = "Agree/Accept" },
    @{ Row = 118; I = "ba"; J = "Appreciation" },
    @{ Row = 119; I = "ba"; J = "Appreciation" },
    @{ Row = 125; I = "ba"; J = "Appreciation" },
    @{ Row = 127; I = "%";  J = "Uninterpretable" },
    @{ Row = 141; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 166; I = "b";  J = "Acknowledge (Backchannel)" },
    @{ Row = 169; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 172; I = "b";  J = "Acknowledge (Backchannel)" },
    @{ Row = 187; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 192; I = "sv"; J = "Statement-opinion" },
    @{ Row = 198; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 221; I = "%";  J = "Uninterpretable" },
    @{ Row = 235; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 241; I = "sv"; J = "Statement-opinion" },
    @{ Row = 242; I = "sv"; J = "Statement-opinion" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}
